# Generate Report for Handback
# Updates the handoff/handback timestamps for the
# "54f9528b-44ec-45a5-9ec9-d8b507c5d970.md" entry across the
# Overview, zh-cn and de-de sheets, as produced by a report-generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
# Row 3 corresponds to 54f9528b-44ec-45a5-9ec9-d8b507c5d970.md
$wsOverview.Range("G3").Value = "2016-11-10 06:40:08"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Row 3 corresponds to 54f9528b-44ec-45a5-9ec9-d8b507c5d970.md
$wsZhCn.Range("H3").Value = "2016-11-10 06:39:53"
$wsZhCn.Range("K3").Value = "2016-11-10 06:40:46"

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
# Row 3 corresponds to 54f9528b-44ec-45a5-9ec9-d8b507c5d970.md
$wsDeDe.Range("H3").Value = "2016-11-10 06:40:08"
$wsDeDe.Range("K3").Value = "2016-11-10 06:41:07"
